$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = '67.335.88'
$ws.Range("E2").Value = '  +0.95%  '

# Row 3
$ws.Range("D3").Value = '3.480.28'
$ws.Range("E3").Value = '  -0.17%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.87'
$ws.Range("E5").Value = '  +0.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.57'
$ws.Range("E6").Value = '  +4.09%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("E8").Value = '  +0.62%  '

# Row 9
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '3.480.26'
$ws.Range("E9").Value = '  -0.16%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  +5.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.08'
$ws.Range("E11").Value = '  -2.27%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.433'
$ws.Range("E12").Value = '  +0.66%  '

# Row 13
$ws.Range("D13").Value = '4.077.84'
$ws.Range("E13").Value = '  -0.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.03'
$ws.Range("E14").Value = '  +11.13%  '

# Row 15
$ws.Range("E15").Value = '  +1.60%  '

# Row 16
$ws.Range("D16").Value = '67.333.91'
$ws.Range("E16").Value = '  +0.89%  '

# Row 17
$ws.Range("E17").Value = '  -0.08%  '

# Row 18
$ws.Range("D18").Value = '3.475.34'
$ws.Range("E18").Value = '  -0.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("E19").Value = '  +0.00%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.29'
$ws.Range("E20").Value = '  +1.80%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.25'
$ws.Range("E21").Value = '  -0.60%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.89'
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.92'
$ws.Range("E23").Value = '  +1.66%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.11%  '

# Row 25
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.535'
$ws.Range("E25").Value = '  +0.28%  '

# Row 26
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.72'
$ws.Range("E26").Value = '  +0.38%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  +0.91%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.36'
$ws.Range("E28").Value = '  +2.04%  '

# Row 29
$ws.Range("E29").Value = '  -3.18%  '

# Row 30
$ws.Range("E30").Value = '  +0.07%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.17'
$ws.Range("E31").Value = '  -0.15%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.42'
$ws.Range("E32").Value = '  -0.14%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.06'
$ws.Range("E33").Value = '  +0.62%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.53'
$ws.Range("E34").Value = '  -0.51%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.36'
$ws.Range("E35").Value = '  +0.58%  '

# Row 36
$ws.Range("E36").Value = '  +0.05%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  -1.60%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.84'
$ws.Range("E38").Value = '  +0.52%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.871'
$ws.Range("E39").Value = '  -0.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  -0.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.73'
$ws.Range("E41").Value = '  +7.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.86'
$ws.Range("E42").Value = '  -0.19%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.63'
$ws.Range("E43").Value = '  -0.30%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.834.41'
$ws.Range("E44").Value = '  +1.45%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.25'
$ws.Range("E45").Value = '  +0.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.84'
$ws.Range("E46").Value = '  -1.08%  '

# Row 47
$ws.Range("E47").Value = '  -2.58%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.57'
$ws.Range("E48").Value = '  -2.52%  '

# Row 49
$ws.Range("E49").Value = '  -1.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '336.30'
$ws.Range("E50").Value = '  +0.28%  '

# Row 51
$ws.Range("E51").Value = '  -1.95%  '
